$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers - copy the existing header style (bold, bordered, centered) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows (row -> I value, J value)
$data = @{
    2  = @(6, 6)
    3  = @(1, 2)
    4  = @(6, 8)
    5  = @(9, 9)
    6  = @(9, 9)
    7  = @(8, 9)
    8  = @(8, 8)
    9  = @(12, 12)
    10 = @(8, 8)
    11 = @(8, 8)
    12 = @(9, 9)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(8, 8)
    16 = @(5, 6)
    17 = @(6, 6)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
